$d = $word.ActiveDocument

$d.Content.Find.Execute("2+39=41", $true, $false, $false, $false, $false, $true, 1, $false, "66+2=68", 2) | Out-Null
$d.Content.Find.Execute("81-14=67", $true, $false, $false, $false, $false, $true, 1, $false, "36+0=36", 2) | Out-Null
$d.Content.Find.Execute("61+29=90", $true, $false, $false, $false, $false, $true, 1, $false, "42-0=42", 2) | Out-Null
$d.Content.Find.Execute("26+31=57", $true, $false, $false, $false, $false, $true, 1, $false, "94-35=59", 2) | Out-Null
$d.Content.Find.Execute("15-7=8", $true, $false, $false, $false, $false, $true, 1, $false, "48-24=24", 2) | Out-Null
$d.Content.Find.Execute("62-5=57", $true, $false, $false, $false, $false, $true, 1, $false, "15-1=14", 2) | Out-Null
$d.Content.Find.Execute("70-14=56", $true, $false, $false, $false, $false, $true, 1, $false, "19+37=56", 2) | Out-Null
$d.Content.Find.Execute("62+12=74", $true, $false, $false, $false, $false, $true, 1, $false, "43-1=42", 2) | Out-Null
$d.Content.Find.Execute("84+10=94", $true, $false, $false, $false, $false, $true, 1, $false, "44+9=53", 2) | Out-Null
$d.Content.Find.Execute("31-18=13", $true, $false, $false, $false, $false, $true, 1, $false, "45+34=79", 2) | Out-Null
$d.Content.Find.Execute("25-23=2", $true, $false, $false, $false, $false, $true, 1, $false, "93-74=19", 2) | Out-Null
$d.Content.Find.Execute("6+83=89", $true, $false, $false, $false, $false, $true, 1, $false, "26+62=88", 2) | Out-Null
$d.Content.Find.Execute("88-69=19", $true, $false, $false, $false, $false, $true, 1, $false, "68+14=82", 2) | Out-Null
$d.Content.Find.Execute("53+11=64", $true, $false, $false, $false, $false, $true, 1, $false, "24+7=31", 2) | Out-Null
$d.Content.Find.Execute("69+12=81", $true, $false, $false, $false, $false, $true, 1, $false, "46-15=31", 2) | Out-Null
$d.Content.Find.Execute("7+86=93", $true, $false, $false, $false, $false, $true, 1, $false, "3+8=11", 2) | Out-Null
$d.Content.Find.Execute("59+4=63", $true, $false, $false, $false, $false, $true, 1, $false, "38+41=79", 2) | Out-Null
$d.Content.Find.Execute("40+39=79", $true, $false, $false, $false, $false, $true, 1, $false, "32+20=52", 2) | Out-Null
$d.Content.Find.Execute("76-28=48", $true, $false, $false, $false, $false, $true, 1, $false, "82-76=6", 2) | Out-Null
$d.Content.Find.Execute("27+27=54", $true, $false, $false, $false, $false, $true, 1, $false, "93-26=67", 2) | Out-Null
$d.Content.Find.Execute("5+78=83", $true, $false, $false, $false, $false, $true, 1, $false, "49-11=38", 2) | Out-Null
$d.Content.Find.Execute("10+55=65", $true, $false, $false, $false, $false, $true, 1, $false, "73+26=99", 2) | Out-Null
$d.Content.Find.Execute("60+29=89", $true, $false, $false, $false, $false, $true, 1, $false, "81-40=41", 2) | Out-Null
$d.Content.Find.Execute("73-69=4", $true, $false, $false, $false, $false, $true, 1, $false, "41+20=61", 2) | Out-Null
$d.Content.Find.Execute("41-30=11", $true, $false, $false, $false, $false, $true, 1, $false, "72+8=80", 2) | Out-Null
$d.Content.Find.Execute("40+53=93", $true, $false, $false, $false, $false, $true, 1, $false, "36-24=12", 2) | Out-Null
$d.Content.Find.Execute("5+3=8", $true, $false, $false, $false, $false, $true, 1, $false, "63+16=79", 2) | Out-Null
$d.Content.Find.Execute("67-48=19", $true, $false, $false, $false, $false, $true, 1, $false, "1+55=56", 2) | Out-Null
$d.Content.Find.Execute("4+81=85", $true, $false, $false, $false, $false, $true, 1, $false, "92-86=6", 2) | Out-Null
$d.Content.Find.Execute("93-21=72", $true, $false, $false, $false, $false, $true, 1, $false, "74-4=70", 2) | Out-Null
$d.Content.Find.Execute("77+6=83", $true, $false, $false, $false, $false, $true, 1, $false, "6+3=9", 2) | Out-Null
$d.Content.Find.Execute("31-23=8", $true, $false, $false, $false, $false, $true, 1, $false, "2+63=65", 2) | Out-Null
$d.Content.Find.Execute("18+26=44", $true, $false, $false, $false, $false, $true, 1, $false, "72+24=96", 2) | Out-Null
$d.Content.Find.Execute("70-54=16", $true, $false, $false, $false, $false, $true, 1, $false, "86-51=35", 2) | Out-Null
$d.Content.Find.Execute("18+60=78", $true, $false, $false, $false, $false, $true, 1, $false, "16+43=59", 2) | Out-Null
$d.Content.Find.Execute("82-62=20", $true, $false, $false, $false, $false, $true, 1, $false, "80-63=17", 2) | Out-Null
$d.Content.Find.Execute("56-46=10", $true, $false, $false, $false, $false, $true, 1, $false, "76+8=84", 2) | Out-Null
$d.Content.Find.Execute("90-66=24", $true, $false, $false, $false, $false, $true, 1, $false, "57+8=65", 2) | Out-Null
$d.Content.Find.Execute("58+30=88", $true, $false, $false, $false, $false, $true, 1, $false, "74-25=49", 2) | Out-Null
$d.Content.Find.Execute("67+17=84", $true, $false, $false, $false, $false, $true, 1, $false, "85-59=26", 2) | Out-Null
$d.Content.Find.Execute("71+22=93", $true, $false, $false, $false, $false, $true, 1, $false, "3+86=89", 2) | Out-Null
$d.Content.Find.Execute("18+39=57", $true, $false, $false, $false, $false, $true, 1, $false, "34+31=65", 2) | Out-Null
$d.Content.Find.Execute("38+4=42", $true, $false, $false, $false, $false, $true, 1, $false, "76-12=64", 2) | Out-Null
$d.Content.Find.Execute("64-29=35", $true, $false, $false, $false, $false, $true, 1, $false, "30-21=9", 2) | Out-Null
$d.Content.Find.Execute("36-28=8", $true, $false, $false, $false, $false, $true, 1, $false, "62-45=17", 2) | Out-Null
$d.Content.Find.Execute("82-3=79", $true, $false, $false, $false, $false, $true, 1, $false, "91-55=36", 2) | Out-Null
$d.Content.Find.Execute("43-19=24", $true, $false, $false, $false, $false, $true, 1, $false, "10+11=21", 2) | Out-Null
$d.Content.Find.Execute("59+20=79", $true, $false, $false, $false, $false, $true, 1, $false, "43+37=80", 2) | Out-Null
$d.Content.Find.Execute("91+3=94", $true, $false, $false, $false, $false, $true, 1, $false, "26+48=74", 2) | Out-Null
$d.Content.Find.Execute("35+5=40", $true, $false, $false, $false, $false, $true, 1, $false, "96-71=25", 2) | Out-Null
$d.Content.Find.Execute("29+20=49", $true, $false, $false, $false, $false, $true, 1, $false, "12+4=16", 2) | Out-Null
$d.Content.Find.Execute("71-21=50", $true, $false, $false, $false, $false, $true, 1, $false, "37-34=3", 2) | Out-Null
$d.Content.Find.Execute("56+4=60", $true, $false, $false, $false, $false, $true, 1, $false, "0+56=56", 2) | Out-Null
$d.Content.Find.Execute("99-31=68", $true, $false, $false, $false, $false, $true, 1, $false, "70+15=85", 2) | Out-Null
$d.Content.Find.Execute("67-49=18", $true, $false, $false, $false, $false, $true, 1, $false, "37-16=21", 2) | Out-Null
$d.Content.Find.Execute("37+6=43", $true, $false, $false, $false, $false, $true, 1, $false, "39+21=60", 2) | Out-Null
$d.Content.Find.Execute("44+42=86", $true, $false, $false, $false, $false, $true, 1, $false, "0+17=17", 2) | Out-Null
$d.Content.Find.Execute("40+29=69", $true, $false, $false, $false, $false, $true, 1, $false, "9+38=47", 2) | Out-Null
$d.Content.Find.Execute("16-3=13", $true, $false, $false, $false, $false, $true, 1, $false, "48+14=62", 2) | Out-Null
$d.Content.Find.Execute("95-7=88", $true, $false, $false, $false, $false, $true, 1, $false, "31+47=78", 2) | Out-Null
$d.Content.Find.Execute("95-72=23", $true, $false, $false, $false, $false, $true, 1, $false, "75-39=36", 2) | Out-Null
$d.Content.Find.Execute("76-33=43", $true, $false, $false, $false, $false, $true, 1, $false, "65-2=63", 2) | Out-Null
$d.Content.Find.Execute("66-32=34", $true, $false, $false, $false, $false, $true, 1, $false, "86-39=47", 2) | Out-Null
$d.Content.Find.Execute("84-13=71", $true, $false, $false, $false, $false, $true, 1, $false, "14+66=80", 2) | Out-Null
$d.Content.Find.Execute("78-29=49", $true, $false, $false, $false, $false, $true, 1, $false, "17+53=70", 2) | Out-Null
$d.Content.Find.Execute("19-6=13", $true, $false, $false, $false, $false, $true, 1, $false, "36-22=14", 2) | Out-Null
$d.Content.Find.Execute("72-27=45", $true, $false, $false, $false, $false, $true, 1, $false, "58+20=78", 2) | Out-Null
$d.Content.Find.Execute("64-12=52", $true, $false, $false, $false, $false, $true, 1, $false, "36+42=78", 2) | Out-Null
$d.Content.Find.Execute("54+24=78", $true, $false, $false, $false, $false, $true, 1, $false, "43+23=66", 2) | Out-Null
$d.Content.Find.Execute("50-41=9", $true, $false, $false, $false, $false, $true, 1, $false, "56-4=52", 2) | Out-Null
$d.Content.Find.Execute("15-0=15", $true, $false, $false, $false, $false, $true, 1, $false, "18+13=31", 2) | Out-Null
$d.Content.Find.Execute("5+83=88", $true, $false, $false, $false, $false, $true, 1, $false, "32+10=42", 2) | Out-Null
$d.Content.Find.Execute("65-34=31", $true, $false, $false, $false, $false, $true, 1, $false, "69-36=33", 2) | Out-Null
$d.Content.Find.Execute("1+43=44", $true, $false, $false, $false, $false, $true, 1, $false, "32-14=18", 2) | Out-Null
$d.Content.Find.Execute("1+62=63", $true, $false, $false, $false, $false, $true, 1, $false, "43+23=66", 2) | Out-Null
$d.Content.Find.Execute("90-52=38", $true, $false, $false, $false, $false, $true, 1, $false, "35-6=29", 2) | Out-Null
$d.Content.Find.Execute("28-5=23", $true, $false, $false, $false, $false, $true, 1, $false, "14+51=65", 2) | Out-Null
$d.Content.Find.Execute("15+53=68", $true, $false, $false, $false, $false, $true, 1, $false, "41-12=29", 2) | Out-Null
$d.Content.Find.Execute("21+21=42", $true, $false, $false, $false, $false, $true, 1, $false, "34-25=9", 2) | Out-Null
$d.Content.Find.Execute("44-30=14", $true, $false, $false, $false, $false, $true, 1, $false, "99-49=50", 2) | Out-Null
$d.Content.Find.Execute("49+1=50", $true, $false, $false, $false, $false, $true, 1, $false, "49-14=35", 2) | Out-Null
$d.Content.Find.Execute("65-28=37", $true, $false, $false, $false, $false, $true, 1, $false, "54-11=43", 2) | Out-Null
$d.Content.Find.Execute("28+16=44", $true, $false, $false, $false, $false, $true, 1, $false, "70-23=47", 2) | Out-Null
$d.Content.Find.Execute("87-44=43", $true, $false, $false, $false, $false, $true, 1, $false, "45+37=82", 2) | Out-Null
$d.Content.Find.Execute("45+39=84", $true, $false, $false, $false, $false, $true, 1, $false, "77+0=77", 2) | Out-Null
$d.Content.Find.Execute("98-28=70", $true, $false, $false, $false, $false, $true, 1, $false, "50+2=52", 2) | Out-Null
$d.Content.Find.Execute("86-48=38", $true, $false, $false, $false, $false, $true, 1, $false, "46-9=37", 2) | Out-Null
$d.Content.Find.Execute("71-34=37", $true, $false, $false, $false, $false, $true, 1, $false, "27+16=43", 2) | Out-Null
$d.Content.Find.Execute("84-57=27", $true, $false, $false, $false, $false, $true, 1, $false, "76+10=86", 2) | Out-Null
$d.Content.Find.Execute("1+52=53", $true, $false, $false, $false, $false, $true, 1, $false, "2+85=87", 2) | Out-Null
$d.Content.Find.Execute("33+4=37", $true, $false, $false, $false, $false, $true, 1, $false, "60-42=18", 2) | Out-Null
$d.Content.Find.Execute("1+57=58", $true, $false, $false, $false, $false, $true, 1, $false, "58-57=1", 2) | Out-Null
$d.Content.Find.Execute("23+39=62", $true, $false, $false, $false, $false, $true, 1, $false, "52-2=50", 2) | Out-Null
$d.Content.Find.Execute("86-58=28", $true, $false, $false, $false, $false, $true, 1, $false, "65-25=40", 2) | Out-Null
$d.Content.Find.Execute("83-38=45", $true, $false, $false, $false, $false, $true, 1, $false, "48-10=38", 2) | Out-Null
$d.Content.Find.Execute("97-8=89", $true, $false, $false, $false, $false, $true, 1, $false, "48+44=92", 2) | Out-Null
$d.Content.Find.Execute("11+45=56", $true, $false, $false, $false, $false, $true, 1, $false, "19+42=61", 2) | Out-Null
$d.Content.Find.Execute("81-58=23", $true, $false, $false, $false, $false, $true, 1, $false, "8+20=28", 2) | Out-Null
$d.Content.Find.Execute("43+49=92", $true, $false, $false, $false, $false, $true, 1, $false, "97-64=33", 2) | Out-Null
$d.Content.Find.Execute("71+6=77", $true, $false, $false, $false, $false, $true, 1, $false, "63+19=82", 2) | Out-Null
